# "Added analysis for close types"
# Adds a second small table (topic/effect pairs) below the existing
# females/males means table, plus a third one-row-header block further
# down, and repositions the existing chart to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First block: topic / effect header + 9 data rows (rows 6-15) ---
$ws.Range("C6").Value = "topic"
$ws.Range("D6").Value = "effect"

$block1 = @(
    @(1, 12.15),
    @(3, 51.75),
    @(4, 111.76),
    @(5, 15.58),
    @(6, 29.14),
    @(8, 12),
    @(9, 88.96),
    @(10, 41.43),
    @(13, 32.07)
)

$r = 7
foreach ($row in $block1) {
    $ws.Cells.Item($r, 3).Value = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $r = $r + 1
}

# --- Second block: marginal header + 2 data rows (rows 18-20) ---
$ws.Range("C18").Value = "marginal"

$block2 = @(
    @(2, 11.13),
    @(12, 7.03)
)

$r = 19
foreach ($row in $block2) {
    $ws.Cells.Item($r, 3).Value = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $r = $r + 1
}

# --- View state: zoom + scroll position + selection ---
$excel.ActiveWindow.Zoom = 90
$ws.Range("B9").Select()

# --- Move/resize the existing chart to make room for the new tables ---
$co = $ws.ChartObjects().Item(1)
$co.Left = 278.52779527559056
$co.Top = 17.25
$co.Width = 952.9375
$co.Height = 394.5
